# Auto-generated edit script: applies cell-value updates described by the
# commit diff across the ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets of the workbook.
# Values are static numeric literals (no formulas in the source sheets), so
# we just overwrite each changed cell's Value.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 11617.25
$ws.Range("J58").Value = 90000
$ws.Range("L58").Value = 270000
$ws.Range("N58").Value = -270300
$ws.Range("H62").Value = 2976.3333
$ws.Range("I62").Value = 3134.1667
$ws.Range("J62").Value = 2660.6667
$ws.Range("K62").Value = 3134.1667
$ws.Range("L62").Value = 2660.6667
$ws.Range("M62").Value = -2510.1667
$ws.Range("N62").Value = -3908.6667
$ws.Range("H64").Value = 2071.4285
$ws.Range("I64").Value = 2000
$ws.Range("J64").Value = 2500
$ws.Range("K64").Value = 2000
$ws.Range("L64").Value = 2500
$ws.Range("M64").Value = -1752
$ws.Range("N64").Value = -2996
$ws.Range("H65").Value = 2976.3333
$ws.Range("I65").Value = 3134.1667
$ws.Range("J65").Value = 2660.6667
$ws.Range("K65").Value = 15670.8335
$ws.Range("L65").Value = 13303.3335
$ws.Range("M65").Value = -12550.8335
$ws.Range("N65").Value = -19543.3335
$ws.Range("H67").Value = 2071.4285
$ws.Range("I67").Value = 2000
$ws.Range("J67").Value = 2500
$ws.Range("K67").Value = 2000
$ws.Range("L67").Value = 2500
$ws.Range("M67").Value = -1142
$ws.Range("N67").Value = -4216
$ws.Range("H76").Value = 44425.125
$ws.Range("I76").Value = 44425.125
$ws.Range("K76").Value = 44425.125
$ws.Range("M76").Value = -44110.125
$ws.Range("H79").Value = 44425.125
$ws.Range("I79").Value = 44425.125
$ws.Range("K79").Value = 44425.125
$ws.Range("M79").Value = -43333.125
$ws.Range("H116").Value = 3166.25
$ws.Range("I116").Value = 1797.5
$ws.Range("J116").Value = 3440
$ws.Range("K116").Value = 1797.5
$ws.Range("L116").Value = 3440
$ws.Range("M116").Value = 1644.5
$ws.Range("N116").Value = -10324
$ws.Range("H125").Value = 1580
$ws.Range("I125").Value = 950
$ws.Range("J125").Value = 2000
$ws.Range("K125").Value = 8550
$ws.Range("L125").Value = 18000
$ws.Range("M125").Value = -6090
$ws.Range("N125").Value = -22920
$ws.Range("H129").Value = 7460.5625
$ws.Range("J129").Value = 11584.2
$ws.Range("L129").Value = 34752.60000000001
$ws.Range("N129").Value = -44752.60000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 48213.97
$ws.Range("I32").Value = 47059.43
$ws.Range("J32").Value = 58989.668
$ws.Range("K32").Value = 47059.43
$ws.Range("L32").Value = 58989.668
$ws.Range("M32").Value = -46772.43
$ws.Range("N32").Value = -59563.668
$ws.Range("H44").Value = 21933.334
$ws.Range("J44").Value = 21933.334
$ws.Range("L44").Value = 21933.334
$ws.Range("N44").Value = -22909.334
$ws.Range("H55").Value = 23800
$ws.Range("J55").Value = 23800
$ws.Range("L55").Value = 23800
$ws.Range("N55").Value = -24430
$ws.Range("H80").Value = 28000
$ws.Range("J80").Value = 28000
$ws.Range("L80").Value = 28000
$ws.Range("N80").Value = -29996
$ws.Range("H83").Value = 28000
$ws.Range("J83").Value = 28000
$ws.Range("L83").Value = 84000
$ws.Range("N83").Value = -93984
$ws.Range("H88").Value = 1117771.8
$ws.Range("I88").Value = 2004867.8
$ws.Range("J88").Value = 8901.75
$ws.Range("K88").Value = 2004867.8
$ws.Range("L88").Value = 8901.75
$ws.Range("M88").Value = -2004461.8
$ws.Range("N88").Value = -9713.75
$ws.Range("H91").Value = 1117771.8
$ws.Range("I91").Value = 2004867.8
$ws.Range("J91").Value = 8901.75
$ws.Range("K91").Value = 2004867.8
$ws.Range("L91").Value = 8901.75
$ws.Range("M91").Value = -2003463.8
$ws.Range("N91").Value = -11709.75
$ws.Range("H132").Value = 6041.8545
$ws.Range("J132").Value = 3000.65
$ws.Range("L132").Value = 9001.950000000001
$ws.Range("N132").Value = -14061.95

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2435.7827
$ws.Range("I86").Value = 3109.9
$ws.Range("J86").Value = 1917.2307
$ws.Range("K86").Value = 3109.9
$ws.Range("L86").Value = 1917.2307
$ws.Range("M86").Value = -1986.9
$ws.Range("N86").Value = -4163.2307
$ws.Range("H89").Value = 2435.7827
$ws.Range("I89").Value = 3109.9
$ws.Range("J89").Value = 1917.2307
$ws.Range("K89").Value = 15549.5
$ws.Range("L89").Value = 9586.1535
$ws.Range("M89").Value = -9933.5
$ws.Range("N89").Value = -20818.1535
$ws.Range("H141").Value = 100000
$ws.Range("J141").Value = 100000
$ws.Range("L141").Value = 100000
$ws.Range("N141").Value = -110360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 11156.571
$ws.Range("I36").Value = 2048
$ws.Range("J36").Value = 14800
$ws.Range("K36").Value = 2048
$ws.Range("L36").Value = 14800
$ws.Range("M36").Value = -1660
$ws.Range("N36").Value = -15576
$ws.Range("H40").Value = 11156.571
$ws.Range("I40").Value = 2048
$ws.Range("J40").Value = 14800
$ws.Range("K40").Value = 2048
$ws.Range("L40").Value = 14800
$ws.Range("M40").Value = -1888
$ws.Range("N40").Value = -15120
$ws.Range("H62").Value = 55559508
$ws.Range("I62").Value = 4350
$ws.Range("J62").Value = 83337090
$ws.Range("K62").Value = 4350
$ws.Range("L62").Value = 83337090
$ws.Range("M62").Value = -3726
$ws.Range("N62").Value = -83338338
$ws.Range("H65").Value = 55559508
$ws.Range("I65").Value = 4350
$ws.Range("J65").Value = 83337090
$ws.Range("K65").Value = 21750
$ws.Range("L65").Value = 416685450
$ws.Range("M65").Value = -18630
$ws.Range("N65").Value = -416691690
$ws.Range("H86").Value = 71430370
$ws.Range("I86").Value = 125001704
$ws.Range("J86").Value = 1916.6666
$ws.Range("K86").Value = 125001704
$ws.Range("L86").Value = 1916.6666
$ws.Range("M86").Value = -125000581
$ws.Range("N86").Value = -4162.6666
$ws.Range("H89").Value = 71430370
$ws.Range("I89").Value = 125001704
$ws.Range("J89").Value = 1916.6666
$ws.Range("K89").Value = 625008520
$ws.Range("L89").Value = 9583.333000000001
$ws.Range("M89").Value = -625002904
$ws.Range("N89").Value = -20815.333
$ws.Range("H123").Value = 61800
$ws.Range("J123").Value = 61800
$ws.Range("L123").Value = 61800
$ws.Range("N123").Value = -71600

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 2143.3333
$ws.Range("I54").Value = 1225
$ws.Range("J54").Value = 3980
$ws.Range("K54").Value = 3675
$ws.Range("L54").Value = 11940
$ws.Range("M54").Value = -3116
$ws.Range("N54").Value = -13058
$ws.Range("H109").Value = 2307.5
$ws.Range("I109").Value = 1880
$ws.Range("K109").Value = 5640
$ws.Range("M109").Value = -4600
$ws.Range("H113").Value = 699.225
$ws.Range("I113").Value = 602.1818
$ws.Range("J113").Value = 1156.7142
$ws.Range("K113").Value = 1806.5454
$ws.Range("L113").Value = 3470.1426
$ws.Range("M113").Value = 363.4546
$ws.Range("N113").Value = -7810.142599999999
$ws.Range("H131").Value = 710.2
$ws.Range("I131").Value = 298.57144
$ws.Range("J131").Value = 777.2093
$ws.Range("K131").Value = 895.71432
$ws.Range("L131").Value = 2331.6279
$ws.Range("M131").Value = 4144.28568
$ws.Range("N131").Value = -12411.6279

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13788091
$ws.Range("I70").Value = 17004242
$ws.Range("J70").Value = 4585.2856
$ws.Range("K70").Value = 17004242
$ws.Range("L70").Value = 4585.2856
$ws.Range("M70").Value = -17003972
$ws.Range("N70").Value = -5125.2856
$ws.Range("H73").Value = 13788091
$ws.Range("I73").Value = 17004242
$ws.Range("J73").Value = 4585.2856
$ws.Range("K73").Value = 17004242
$ws.Range("L73").Value = 4585.2856
$ws.Range("M73").Value = -17003306
$ws.Range("N73").Value = -6457.2856
$ws.Range("H80").Value = 12000
$ws.Range("I80").Value = 9000
$ws.Range("J80").Value = 18000
$ws.Range("K80").Value = 9000
$ws.Range("L80").Value = 18000
$ws.Range("M80").Value = -8002
$ws.Range("N80").Value = -19996
$ws.Range("H83").Value = 12000
$ws.Range("I83").Value = 9000
$ws.Range("J83").Value = 18000
$ws.Range("K83").Value = 45000
$ws.Range("L83").Value = 90000
$ws.Range("M83").Value = -40008
$ws.Range("N83").Value = -99984
$ws.Range("H97").Value = 1408.5366
$ws.Range("I97").Value = 1226.3334
$ws.Range("J97").Value = 1905.4546
$ws.Range("K97").Value = 1226.3334
$ws.Range("L97").Value = 1905.4546
$ws.Range("M97").Value = -730.3334
$ws.Range("N97").Value = -2897.4546
$ws.Range("H123").Value = 16118.2
$ws.Range("J123").Value = 16118.2
$ws.Range("L123").Value = 16118.2
$ws.Range("N123").Value = -21018.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6767.4653
$ws.Range("I132").Value = 7505.8
$ws.Range("J132").Value = 3537.25
$ws.Range("K132").Value = 22517.4
$ws.Range("L132").Value = 10611.75
$ws.Range("M132").Value = -19987.4
$ws.Range("N132").Value = -15671.75

